# Update cryptos list data (prices and 1h volume changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.234.65"
$ws.Range('E2').Value = '  -0.97%  '

$ws.Range('D3').Value = "'1.860.35"
$ws.Range('E3').Value = '  -0.65%  '

$ws.Range('D4').Value = "'0.9999"
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').Value = "'0.7152"
$ws.Range('E5').Value = '  -0.55%  '

$ws.Range('D6').Value = "'240.56"
$ws.Range('E6').Value = '  +0.52%  '

$ws.Range('D8').Value = "'0.07753"
$ws.Range('E8').Value = '  -1.15%  '

$ws.Range('D9').Value = "'0.3081"
$ws.Range('E9').Value = '  +0.25%  '

$ws.Range('D10').Value = "'25.20"
$ws.Range('E10').Value = '  -0.60%  '

$ws.Range('D11').Value = "'0.08259"
$ws.Range('E11').Value = '  +0.32%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = "'1.859.97"
$ws.Range('E12').Value = '  -1.17%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'5.240"
$ws.Range('E13').Value = '  +0.02%  '

$ws.Range('D14').Value = "'0.7182"
$ws.Range('E14').Value = '  -0.65%  '

$ws.Range('D15').Value = "'90.30"
$ws.Range('E15').Value = '  +0.44%  '

$ws.Range('D16').Value = "'29.219.36"
$ws.Range('E16').Value = '  -1.45%  '

$ws.Range('D17').Value = "'5.867"
$ws.Range('E17').Value = '  +0.63%  '

$ws.Range('D18').Value = "'244.39"
$ws.Range('E18').Value = '  +1.15%  '

$ws.Range('D19').Value = "'0.000007800"
$ws.Range('E19').Value = '  -0.71%  '

$ws.Range('D20').Value = "'13.16"
$ws.Range('E20').Value = '  -1.27%  '

$ws.Range('D21').Value = "'2.108.51"
$ws.Range('E21').Value = '  -1.01%  '

$ws.Range('D22').Value = "'0.9999"
$ws.Range('E22').Value = '  -0.04%  '

$ws.Range('D23').Value = "'7.984"
$ws.Range('E23').Value = '  +3.21%  '

$ws.Range('D24').Value = "'1.000"
$ws.Range('E24').Value = '  -0.12%  '

$ws.Range('E25').Value = '  +1.74%  '

$ws.Range('D26').Value = "'162.49"
$ws.Range('E26').Value = '  -0.07%  '

$ws.Range('E27').Value = '  -0.12%  '

$ws.Range('D28').Value = "'18.26"
$ws.Range('E28').Value = '  -0.31%  '

$ws.Range('E29').Value = '  +0.86%  '

$ws.Range('D30').Value = "'1.314"
$ws.Range('E30').Value = '  -3.47%  '

$ws.Range('D31').Value = "'4.405"
$ws.Range('E31').Value = '  +1.51%  '

$ws.Range('D32').Value = "'4.184"
$ws.Range('E32').Value = '  +2.67%  '

$ws.Range('D33').Value = "'0.05193"
$ws.Range('E33').Value = '  -1.23%  '

$ws.Range('E34').Value = '  -1.34%  '

$ws.Range('D35').Value = "'1.172"
$ws.Range('E35').Value = '  -2.10%  '

$ws.Range('D36').Value = "'0.7281"
$ws.Range('E36').Value = '  +1.65%  '

$ws.Range('D37').Value = "'2.677"
$ws.Range('E37').Value = '  +0.29%  '

$ws.Range('E38').Value = '  -0.68%  '

$ws.Range('D39').Value = "'2.687"
$ws.Range('E39').Value = '  -1.27%  '

$ws.Range('D40').Value = "'1.150.41"
$ws.Range('E40').Value = '  -2.43%  '

$ws.Range('D41').Value = "'0.9055"
$ws.Range('E41').Value = '  -0.32%  '

$ws.Range('D42').Value = "'6.094"
$ws.Range('E42').Value = '  +1.59%  '

$ws.Range('D43').Value = "'72.35"
$ws.Range('E43').Value = '  +1.12%  '

$ws.Range('D44').Value = "'1.0000"
$ws.Range('E44').Value = '  -0.12%  '

$ws.Range('D45').Value = "'101.79"
$ws.Range('E45').Value = '  -0.71%  '

$ws.Range('D46').Value = "'2.008.33"
$ws.Range('E46').Value = '  -1.00%  '

$ws.Range('D47').Value = "'0.5231"
$ws.Range('E47').Value = '  -2.41%  '

$ws.Range('E48').Value = '  +0.13%  '

$ws.Range('D49').Value = "'9.334"
$ws.Range('E49').Value = '  +1.82%  '

$ws.Range('D50').Value = "'2.873"
$ws.Range('E50').Value = '  +1.23%  '

$ws.Range('D51').Value = "'7.060"
$ws.Range('E51').Value = '  +0.46%  '
